$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "data group" (columns A-R) down by one row for rows 91..169,
# working from the bottom up so we don't clobber values before they're copied.
for ($r = 168; $r -ge 90; $r--) {
    $src = $r
    $dst = $r + 1

    $ws.Cells.Item($dst, 4).Value  = $ws.Cells.Item($src, 4).Value2()   # D - Fecha
    $ws.Cells.Item($dst, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat()
    $ws.Cells.Item($dst, 10).Value = $ws.Cells.Item($src, 10).Value2()  # J - Volumen
    $ws.Cells.Item($dst, 11).Value = $ws.Cells.Item($src, 11).Value2()  # K - Precio minimo
    $ws.Cells.Item($dst, 12).Value = $ws.Cells.Item($src, 12).Value2()  # L - Precio maximo
    $ws.Cells.Item($dst, 13).Value = $ws.Cells.Item($src, 13).Value2()  # M - Precio promedio ponderado
    $ws.Cells.Item($dst, 16).Value = $ws.Cells.Item($src, 16).Value2()  # P - Precio $/Kg

    # Columns A, B, C, E, F, G, H, I, N, O, Q, R are identical for every row in this
    # block already, but copy them along too so row 169 is fully populated.
    $ws.Cells.Item($dst, 1).Value  = $ws.Cells.Item($src, 1).Value2()
    $ws.Cells.Item($dst, 2).Value  = $ws.Cells.Item($src, 2).Value2()
    $ws.Cells.Item($dst, 3).Value  = $ws.Cells.Item($src, 3).Value2()
    $ws.Cells.Item($dst, 5).Value  = $ws.Cells.Item($src, 5).Value2()
    $ws.Cells.Item($dst, 6).Value  = $ws.Cells.Item($src, 6).Value2()
    $ws.Cells.Item($dst, 7).Value  = $ws.Cells.Item($src, 7).Value2()
    $ws.Cells.Item($dst, 8).Value  = $ws.Cells.Item($src, 8).Value2()
    $ws.Cells.Item($dst, 9).Value  = $ws.Cells.Item($src, 9).Value2()
    $ws.Cells.Item($dst, 14).Value = $ws.Cells.Item($src, 14).Value2()
    $ws.Cells.Item($dst, 15).Value = $ws.Cells.Item($src, 15).Value2()
    $ws.Cells.Item($dst, 17).Value = $ws.Cells.Item($src, 17).Value2()
    $ws.Cells.Item($dst, 18).Value = $ws.Cells.Item($src, 18).Value2()
}

# Row 90 keeps its other values but gets a brand-new date.
$ws.Cells.Item(90, 4).Value = 44447
